$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.614.12'
$ws.Range('D3').Value = '2.498.87'
$ws.Range('E3').Value = '  +2.90%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '323.84'
$ws.Range('E5').Value = '  +2.09%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '105.23'
$ws.Range('E6').Value = '  +2.44%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.523'
$ws.Range('E7').Value = '  +1.93%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +2.58%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '37.76'
$ws.Range('E10').Value = '  +6.58%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0813'
$ws.Range('E11').Value = '  +1.49%  '
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '18.29'
$ws.Range('E13').Value = '  +0.78%  '
$ws.Range('E14').Value = '  +1.90%  '
$ws.Range('D15').Value = '2.887.86'
$ws.Range('E15').Value = '  +2.84%  '
$ws.Range('D16').Value = '2.497.78'
$ws.Range('E16').Value = '  +2.97%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.843'
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').Value = '47.482.73'
$ws.Range('E18').Value = '  +5.41%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.66'
$ws.Range('E19').Value = '  +3.36%  '
$ws.Range('E20').Value = '  +2.91%  '
$ws.Range('D21').Value = '0.0₃0934'
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '70.80'
$ws.Range('E22').Value = '  +2.83%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '251.22'
$ws.Range('E23').Value = '  +3.04%  '
$ws.Range('E24').Value = '  +5.58%  '
$ws.Range('E25').Value = '  +2.95%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '26.25'
$ws.Range('E26').Value = '  +3.81%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.09'
$ws.Range('E28').Value = '  +5.63%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.20'
$ws.Range('E29').Value = '  +0.67%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '35.28'
$ws.Range('E30').Value = '  +7.78%  '
$ws.Range('E32').Value = '  +0.49%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '19.99'
$ws.Range('E33').Value = '  -1.17%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.38'
$ws.Range('E34').Value = '  +3.21%  '
$ws.Range('E35').Value = '  +2.65%  '
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('E37').Value = '  +4.15%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.64'
$ws.Range('E38').Value = '  +5.03%  '
$ws.Range('E39').Value = '  +4.09%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.112'
$ws.Range('E40').Value = '  +2.31%  '
$ws.Range('E41').Value = '  +0.47%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '121.13'
$ws.Range('E42').Value = '  -3.27%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '21.09'
$ws.Range('E43').Value = '  +1.39%  '
$ws.Range('E44').Value = '  +2.74%  '
$ws.Range('D45').Value = '1.969.16'
$ws.Range('E45').Value = '  +1.83%  '
$ws.Range('E46').Value = '  +1.67%  '
$ws.Range('E47').Value = '  -0.48%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.85'
$ws.Range('E48').Value = '  +1.72%  '
$ws.Range('E49').Value = '  -0.49%  '
$ws.Range('E50').Value = '  +13.48%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '79.09'
$ws.Range('E51').Value = '  +3.47%  '

Write-Output "Applied all cell updates."